# NGIN & VLV updated
# Update the NGIN (8723 -> 1020) and SAN (1223782 -> 1223754) sample data
# values that are repeated across row 2 of the "NGIN" sheet, move the
# selection / scroll position, and resize column C (no longer best-fit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data refresh: NGIN8723 family -> NGIN1020, SAN 1223782 -> 1223754 ---
$ws.Range("A2").Value  = "NGIN1020"
$ws.Range("C2").Value  = "NGIN1020"
$ws.Range("D2").Value  = "ngindomain1020.com"
$ws.Range("F2").Value  = "nginocn1020"
$ws.Range("G2").Value  = "testreference1020"
$ws.Range("H2").Value  = "ngincontact1020"
$ws.Range("J2").Value  = "ngin1020@test.com"
$ws.Range("Q2").Value  = "NGIN1020"
$ws.Range("T2").Value  = "NGINOrder_1020"
$ws.Range("U2").Value  = "NGINRFI_1020"
$ws.Range("V2").Value  = "NGINOrder_1020"
$ws.Range("W2").Value  = "NGINRFI_1020"
$ws.Range("Y2").Value  = "NGINService_1020"
$ws.Range("AJ2").Value = "NGINUser_1020"
$ws.Range("AK2").Value = "User_1020"
$ws.Range("AN2").Value = "NGINUser_1020@gmail.com"
$ws.Range("AP2").Value = "NGINOrder_1020"
$ws.Range("AQ2").Value = "NGINUseredit1020"
$ws.Range("AR2").Value = "Useredit1020"
$ws.Range("AU2").Value = "NGINUseredit_1020@gmail.com"
$ws.Range("AZ2").Value = "NGINOrderedit_1020"
$ws.Range("BA2").Value = "NGINRFIedit_1020"
$ws.Range("BB2").Value = "NGINOrder_1020"
$ws.Range("BC2").Value = "NGINRFI_1020"
$ws.Range("BI2").Value = "nginreseller1020@gmail.com"
$ws.Range("BQ2").Value = "nginreselleredit1020@gmail.com"
$ws.Range("BY2").Value = "Reseller1020"
$ws.Range("BZ2").Value = "1223754"
$ws.Range("CN2").Value = "1223754"
$ws.Range("DP2").Value = "Reselleredit1020"
$ws.Range("ED2").Value = "AT-nginocn1020"
$ws.Range("EE2").Value = "431223754"
$ws.Range("EF2").Value = "390201020891"

# --- Column C is no longer best-fit; it was manually narrowed ---
$ws.Columns("C").ColumnWidth = 8

# --- Move the view: scroll toward column DZ and select EE6 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 130
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("EE6").Select()
